$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (Amplificatore Operazionale): swap the component from TL071 (Texas
# Instruments / SAMB, no ref) to AD826AN (Analog Devices / Distrelec /
# 173-22-437), and update its unit price.
$ws.Range("B12").Value = "AD826AN"
$ws.Range("C12").Value = "Analog Devices"
$ws.Range("D12").Value = "Distrelec"
$ws.Range("E12").Value = "173-22-437"
$ws.Range("G12").Value = 8.9

# Row 19 (Trimmer 25K / 67WR25KLF / BI Technologies) is removed from the
# BOM entirely; deleting the whole row shifts everything below it up by
# one, which also moves the "Totale:" row from 21 to 20 and keeps its
# SUM formula correctly scoped to the new last data row.
$ws.Rows("19").Delete()

# Widen column H slightly.
$ws.Columns(8).ColumnWidth = 9.8

# Restore the cursor to where the editor left it.
$ws.Range("F17").Select()
